$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values differ between row 2 and row 3 and need to be swapped
$cols = @("A", "B", "D", "E", "F", "G", "H", "AH")

foreach ($col in $cols) {
    $r2 = $ws.Range($col + "2")
    $r3 = $ws.Range($col + "3")
    $v2 = $r2.Value2
    $v3 = $r3.Value2
    $r2.Value2 = $v3
    $r3.Value2 = $v2
}
